$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the sheet (the workbook/sheet metadata is being re-pointed at the
#    "Informe-01-010040-TC-TM" cube instead of "Informe-05-050314-A-TC-TP").
# ---------------------------------------------------------------------------
$ws.Name = "Informe-01-010040-TC-TM"

# ---------------------------------------------------------------------------
# 2. The sheet's used area shrinks from 19 data columns (A:S) + 6 unused,
#    formatted-only columns (T:Y) down to just 15 data columns (A:O) plus two
#    blank, still-formatted columns (P:Q maps to old R:S... see below).
#    Concretely: columns A:O get brand-new widths, and the old custom widths
#    that used to live in P:S (old indices 16-19) are dropped back to the
#    sheet's normal/default width.
# ---------------------------------------------------------------------------
$newWidths = @(27.69, 44.5, 18.66, 55.2, 34.64, 36.31, 47.28, 27.69, 27.69, 27.69, 15.46, 46.44, 19.19, 20.05, 29.5)
for ($i = 0; $i -lt $newWidths.Length; $i++) {
    # The engine stores column widths on a pixel grid (Excel's classic
    # "MDW" character-width model), so we compensate for the fixed ~0.8333
    # offset that it always re-adds when persisting the value.
    $ws.Columns.Item($i + 1).ColumnWidth = ($newWidths[$i] - 0.8333333333333333)
}

# Columns P:S (16-19) used to carry bespoke widths; they revert to the plain
# sheet default (~11.52) now that they're outside the "important" range.
for ($i = 16; $i -le 19; $i++) {
    $ws.Columns.Item($i).ColumnWidth = (11.52 - 0.8333333333333333)
}

# ---------------------------------------------------------------------------
# 3. Columns T:Y (20-25) are no longer part of the sheet at all -- the cells
#    that used to sit there (rows 1-5, formatted but empty) are removed.
# ---------------------------------------------------------------------------
$ws.Range("T1:Y5").Clear()

# ---------------------------------------------------------------------------
# 4. The sheet's stored selection/view now spans A1:Q10 instead of just A1.
# ---------------------------------------------------------------------------
$ws.Range("A1:Q10").Select()

# ---------------------------------------------------------------------------
# 5. Six extra (still empty) rows get appended below the existing data
#    (rows 9-14), matching the height already used by the existing blank
#    rows 7-8.
# ---------------------------------------------------------------------------
for ($r = 9; $r -le 14; $r++) {
    $ws.Rows.Item($r).RowHeight = 12.8
}
